$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Te gusta" label to include the rest of the sentence.
$ws.Range("E9").Value2 = "Te gusta esta página"

# Move the active selection from F6 to E10 (last thing the author clicked
# before saving).
$ws.Range("E10").Select() | Out-Null
